# Adds a new "2022-Q4" quarter sheet (right after "总计") with the latest
# fund-holding snapshot, and inserts a matching new row at the top of the
# "总计" (totals) sheet's data table.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Update "总计" sheet: insert a new row 2 for "2022-Q4" and shift the
#    existing quarterly rows down by one (their own data follows them).
# ---------------------------------------------------------------------
$totals = $wb.Worksheets.Item(1)

$totals.Rows(2).Insert()
# Copy the (plain, un-bolded) format of the row below onto the new row 2
# so it doesn't inherit the bold header style that Insert() copies from
# row 1 above.
$totals.Range("A3").Copy()
$totals.Range("A2").PasteSpecial(-4122)
$totals.Range("B2:D2").ClearFormats()

$totals.Range("A2").Value = 0
$totals.Range("B2").Value = "2022-Q4"
$totals.Range("C2").Value = 1
$totals.Range("D2").Value = 1.93

# Renumber the A column (0-based sequential index) for every data row.
for ($r = 2; $r -le 10; $r++) {
    $totals.Cells.Item($r, 1).Value = $r - 2
}

# ---------------------------------------------------------------------
# 2. Insert the new "2022-Q4" worksheet right after "总计".
# ---------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add($null, $totals)
$newSheet.Name = "2022-Q4"

$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

$newSheet.Range("A2").Value = 0

# B2..G2 must be stored as TEXT (matching every other quarter sheet),
# so force text format before assigning, then drop the incidental style
# that NumberFormat="@" leaves behind.
$newSheet.Range("B2:G2").NumberFormat = "@"
$newSheet.Range("B2").Value = "510810"
$newSheet.Range("C2").Value = "汇添富中证上海国企ETF"
$newSheet.Range("D2").Value = "65.67"
$newSheet.Range("E2").Value = "98.06"
$newSheet.Range("F2").Value = "2.94"
$newSheet.Range("G2").Value = "1.9307"
$newSheet.Range("B2:G2").ClearFormats()

$newSheet.Range("H2").Value = 9

# ---------------------------------------------------------------------
# 3. Restore the originally-selected tab ("2020-Q4", the last sheet) as
#    the active tab, since adding a sheet makes it active by default.
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$lastSheet.Activate()
